$changes = @(
  @{ Sheet = "Citywide Totals"; Cells = @(
      @{ Row = 2; NewValue = 4969 },
      @{ Row = 3; NewValue = 5111 },
      @{ Row = 4; NewValue = 1061 },
      @{ Row = 5; NewValue = 361 },
      @{ Row = 6; NewValue = 5740 },
      @{ Row = 7; NewValue = 17242 }
  )},
  @{ Sheet = "Austin"; Cells = @(
      @{ Row = 3; NewValue = 342 },
      @{ Row = 6; NewValue = 393 },
      @{ Row = 7; NewValue = 1152 }
  )},
  @{ Sheet = "South Chicago"; Cells = @(
      @{ Row = 2; NewValue = 128 },
      @{ Row = 3; NewValue = 140 },
      @{ Row = 7; NewValue = 381 }
  )},
  @{ Sheet = "Garfield Park"; Cells = @(
      @{ Row = 2; NewValue = 200 },
      @{ Row = 3; NewValue = 271 },
      @{ Row = 6; NewValue = 211 },
      @{ Row = 7; NewValue = 730 }
  )},
  @{ Sheet = "Grand Crossing"; Cells = @(
      @{ Row = 3; NewValue = 191 },
      @{ Row = 4; NewValue = 25 }
  )},
  @{ Sheet = "New City"; Cells = @(
      @{ Row = 3; NewValue = 99 },
      @{ Row = 4; NewValue = 13 },
      @{ Row = 6; NewValue = 152 },
      @{ Row = 7; NewValue = 393 }
  )},
  @{ Sheet = "Woodlawn"; Cells = @(
      @{ Row = 2; NewValue = 75 },
      @{ Row = 3; NewValue = 121 },
      @{ Row = 7; NewValue = 296 }
  )},
  @{ Sheet = "By Neighborhood"; Cells = @(
      @{ Row = 6; NewValue = 130 },
      @{ Row = 7; NewValue = 515 },
      @{ Row = 8; NewValue = 1152 },
      @{ Row = 14; NewValue = 93 },
      @{ Row = 15; NewValue = 173 },
      @{ Row = 20; NewValue = 394 },
      @{ Row = 27; NewValue = 158 },
      @{ Row = 29; NewValue = 925 },
      @{ Row = 32; NewValue = 23 },
      @{ Row = 33; NewValue = 730 },
      @{ Row = 34; NewValue = 92 },
      @{ Row = 36; NewValue = 225 },
      @{ Row = 42; NewValue = 642 },
      @{ Row = 43; NewValue = 152 },
      @{ Row = 44; NewValue = 152 },
      @{ Row = 48; NewValue = 217 },
      @{ Row = 49; NewValue = 97 },
      @{ Row = 51; NewValue = 218 },
      @{ Row = 54; NewValue = 338 },
      @{ Row = 55; NewValue = 196 },
      @{ Row = 57; NewValue = 63 },
      @{ Row = 63; NewValue = 51 },
      @{ Row = 64; NewValue = 109 },
      @{ Row = 65; NewValue = 393 },
      @{ Row = 66; NewValue = 57 },
      @{ Row = 67; NewValue = 665 },
      @{ Row = 71; NewValue = 55 },
      @{ Row = 73; NewValue = 147 },
      @{ Row = 75; NewValue = 59 },
      @{ Row = 79; NewValue = 420 },
      @{ Row = 83; NewValue = 381 },
      @{ Row = 84; NewValue = 131 },
      @{ Row = 85; NewValue = 793 },
      @{ Row = 86; NewValue = 116 },
      @{ Row = 89; NewValue = 250 },
      @{ Row = 94; NewValue = 230 },
      @{ Row = 96; NewValue = 184 },
      @{ Row = 99; NewValue = 296 },
      @{ Row = 101; NewValue = 17242 }
  )},
  @{ Sheet = "North Lawndale"; Cells = @(
      @{ Row = 3; NewValue = 232 },
      @{ Row = 6; NewValue = 190 },
      @{ Row = 7; NewValue = 665 }
  )},
  @{ Sheet = "South Deering"; Cells = @(
      @{ Row = 3; NewValue = 52 },
      @{ Row = 7; NewValue = 131 }
  )},
  @{ Sheet = "Lincoln Park"; Cells = @(
      @{ Row = 3; NewValue = 19 },
      @{ Row = 7; NewValue = 97 }
  )},
  @{ Sheet = "Loop"; Cells = @(
      @{ Row = 6; NewValue = 179 },
      @{ Row = 7; NewValue = 338 }
  )},
  @{ Sheet = "Englewood"; Cells = @(
      @{ Row = 2; NewValue = 268 },
      @{ Row = 6; NewValue = 257 },
      @{ Row = 7; NewValue = 925 }
  )},
  @{ Sheet = "Lake View"; Cells = @(
      @{ Row = 6; NewValue = 107 },
      @{ Row = 7; NewValue = 217 }
  )},
  @{ Sheet = "Irving Park"; Cells = @(
      @{ Row = 2; NewValue = 37 },
      @{ Row = 7; NewValue = 152 }
  )},
  @{ Sheet = "Bridgeport"; Cells = @(
      @{ Row = 3; NewValue = 20 },
      @{ Row = 7; NewValue = 93 }
  )},
  @{ Sheet = "Ashburn"; Cells = @(
      @{ Row = 2; NewValue = 49 },
      @{ Row = 7; NewValue = 130 }
  )},
  @{ Sheet = "Humboldt Park"; Cells = @(
      @{ Row = 2; NewValue = 168 },
      @{ Row = 3; NewValue = 201 },
      @{ Row = 4; NewValue = 25 },
      @{ Row = 5; NewValue = 5 },
      @{ Row = 7; NewValue = 642 }
  )},
  @{ Sheet = "Lower West Side"; Cells = @(
      @{ Row = 3; NewValue = 57 },
      @{ Row = 7; NewValue = 196 }
  )},
  @{ Sheet = "West Ridge"; Cells = @(
      @{ Row = 6; NewValue = 80 },
      @{ Row = 7; NewValue = 184 }
  )},
  @{ Sheet = "Roseland"; Cells = @(
      @{ Row = 3; NewValue = 135 },
      @{ Row = 7; NewValue = 420 }
  )},
  @{ Sheet = "Near South Side"; Cells = @(
      @{ Row = 6; NewValue = 42 },
      @{ Row = 7; NewValue = 109 }
  )},
  @{ Sheet = "Chicago Lawn"; Cells = @(
      @{ Row = 2; NewValue = 132 },
      @{ Row = 7; NewValue = 394 }
  )},
  @{ Sheet = "Grand Boulevard"; Cells = @(
      @{ Row = 2; NewValue = 90 },
      @{ Row = 7; NewValue = 225 }
  )},
  @{ Sheet = "Auburn Gresham"; Cells = @(
      @{ Row = 2; NewValue = 177 },
      @{ Row = 3; NewValue = 168 },
      @{ Row = 7; NewValue = 515 }
  )},
  @{ Sheet = "Garfield Ridge"; Cells = @(
      @{ Row = 3; NewValue = 25 },
      @{ Row = 7; NewValue = 92 }
  )},
  @{ Sheet = "West Loop"; Cells = @(
      @{ Row = 6; NewValue = 99 },
      @{ Row = 7; NewValue = 230 }
  )},
  @{ Sheet = "Brighton Park"; Cells = @(
      @{ Row = 6; NewValue = 54 },
      @{ Row = 7; NewValue = 173 }
  )},
  @{ Sheet = "North Center"; Cells = @(
      @{ Row = 4; NewValue = 1 },
      @{ Row = 7; NewValue = 57 }
  )},
  @{ Sheet = "Portage Park"; Cells = @(
      @{ Row = 2; NewValue = 47 },
      @{ Row = 7; NewValue = 147 }
  )},
  @{ Sheet = "Galewood"; Cells = @(
      @{ Row = 2; NewValue = 8 },
      @{ Row = 7; NewValue = 23 }
  )},
  @{ Sheet = "Uptown"; Cells = @(
      @{ Row = 4; NewValue = 28 },
      @{ Row = 7; NewValue = 250 }
  )},
  @{ Sheet = "Edgewater"; Cells = @(
      @{ Row = 2; NewValue = 41 },
      @{ Row = 6; NewValue = 60 },
      @{ Row = 7; NewValue = 158 }
  )},
  @{ Sheet = "Streeterville"; Cells = @(
      @{ Row = 6; NewValue = 29 },
      @{ Row = 7; NewValue = 116 }
  )},
  @{ Sheet = "Pullman"; Cells = @(
      @{ Row = 6; NewValue = 10 },
      @{ Row = 7; NewValue = 59 }
  )},
  @{ Sheet = "Little Italy, UIC"; Cells = @(
      @{ Row = 2; NewValue = 63 },
      @{ Row = 6; NewValue = 74 },
      @{ Row = 7; NewValue = 218 }
  )},
  @{ Sheet = "Mckinley Park"; Cells = @(
      @{ Row = 4; NewValue = 4 },
      @{ Row = 7; NewValue = 63 }
  )},
  @{ Sheet = "Hyde Park"; Cells = @(
      @{ Row = 3; NewValue = 42 },
      @{ Row = 7; NewValue = 152 }
  )},
  @{ Sheet = "South Shore"; Cells = @(
      @{ Row = 4; NewValue = 48 },
      @{ Row = 6; NewValue = 190 },
      @{ Row = 7; NewValue = 793 }
  )},
  @{ Sheet = "Oakland"; Cells = @(
      @{ Row = 4; NewValue = 2 },
      @{ Row = 7; NewValue = 55 }
  )}
)

$wb = $excel.ActiveWorkbook

foreach ($entry in $changes) {
    $ws = $wb.Worksheets.Item($entry.Sheet)
    foreach ($cell in $entry.Cells) {
        $ws.Cells.Item($cell.Row, 11).Value = $cell.NewValue
    }
}
